# Update "想去人数" (F column) counts for several conventions/events across
# all four worksheets, reflecting the latest scrape snapshot.

$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1): rows keyed by row number -> new F value
$zhanlan = @{
    2  = 11
    3  = 7626
    4  = 2712
    5  = 930
    6  = 287
    7  = 795
    8  = 582
    9  = 85
    11 = 377
    13 = 3062
    14 = 189
    15 = 81
    16 = 715
    19 = 445
    21 = 203
    22 = 211
    23 = 259
    24 = 283
    25 = 125
    26 = 92
    27 = 243
    30 = 490
    31 = 444
    35 = 82
}

# Sheet "演出" (sheet2)
$yanchu = @{
    7 = 8
}

# Sheet "本地生活" (sheet3)
$bendi = @{
    2 = 203
}

# Sheet "全部类型" (sheet4)
$quanbu = @{
    2  = 203
    5  = 11
    6  = 7626
    7  = 2712
    8  = 930
    9  = 287
    10 = 795
    11 = 582
    12 = 85
    14 = 377
    17 = 3062
    18 = 189
    19 = 81
    21 = 715
    25 = 445
    27 = 203
    28 = 211
    29 = 259
    30 = 283
    31 = 125
    32 = 92
    33 = 243
    36 = 490
    37 = 444
    41 = 82
    42 = 8
}

function Apply-Updates($sheetName, $updates) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}

Apply-Updates "展览" $zhanlan
Apply-Updates "演出" $yanchu
Apply-Updates "本地生活" $bendi
Apply-Updates "全部类型" $quanbu
